# Added url to wireframe workbook:
#  - A1 label changed from "Wireframes" to "Wireframe"
#  - B1 gets a new hyperlink to https://wireframe.cc/Hpuwsw (with the
#    standard built-in "Hyperlink" style: underlined, theme-10 colored font)
#  - Selection cursor moved to C13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value2 = "Wireframe"

[void]$ws.Hyperlinks.Add($ws.Range("B1"), "https://wireframe.cc/Hpuwsw")

[void]$ws.Range("C13").Select()
